# Regenerate sval data (filter save games) - update numeric columns B:G
# for rows 2-12 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.230985683306322, 10.29869402782916, 0.8054896365839992, 8.660232485948974, 0, 22.99540183366846)
    3  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    4  = @(1.459612070389937, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 2.42670696938877)
    5  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797)
    6  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    7  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    8  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    9  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    10 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 6.740334628841572)
    11 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 5.553084769722144)
    12 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
